$wb = $excel.ActiveWorkbook

# --- Sheet names ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "T1 - TEst"
$ws2.Name = "T2 - Test 2"

# --- Sheet1 header / title cells ---
$ws1.Range("A1").Value = "TEst"
$ws1.Range("A2").Value = "Test Frame"
$ws1.Range("F3").Value = "Test"
$ws1.Range("H3").Value = "test 2"

# --- Sheet2 header / title cells ---
$ws2.Range("A1").Value = "Test 2"
$ws2.Range("A2").Value = "Test Frame 2"
$ws2.Range("F3").Value = "Test"
$ws2.Range("H3").Value = "test 2"

# --- Sheet1 column F (Mindset/Test) data updates ---
$ws1.Range("F4").Value = "Impassioned"
$ws1.Range("F5").Value = "Impassioned"
$ws1.Range("F6").Value = "Impassioned"
$ws1.Range("F7").Value = "Attracted"
$ws1.Range("F8").Value = "Attracted"
$ws1.Range("F9").Value = "Impassioned"
$ws1.Range("F10").Value = "Impassioned"
$ws1.Range("F11").Value = "Attracted"
$ws1.Range("F12").Value = "Impassioned"
$ws1.Range("F13").Value = "Attracted"
$ws1.Range("F15").Value = "Impassioned"
$ws1.Range("F16").Value = "Impassioned"
$ws1.Range("F17").Value = "Attracted"
$ws1.Range("F18").Value = "Attracted"
$ws1.Range("F19").Value = "Attracted"
$ws1.Range("F20").Value = "Attracted"
$ws1.Range("F21").Value = "Attracted"
$ws1.Range("F22").Value = "Attracted"
$ws1.Range("F23").Value = "Attracted"
$ws1.Range("F25").Value = "Impassioned"
$ws1.Range("F26").Value = "Impassioned"
$ws1.Range("F28").Value = "Attracted"
$ws1.Range("F29").Value = "Attracted"
$ws1.Range("F31").Value = "Attracted"
$ws1.Range("F32").Value = "Apathetic"
$ws1.Range("F33").Value = "Attracted"
$ws1.Range("F34").Value = "Impassioned"
$ws1.Range("F36").Value = "Impassioned"
$ws1.Range("F37").Value = "Apathetic"
$ws1.Range("F38").Value = "Apathetic"
$ws1.Range("F39").Value = "Apathetic"
$ws1.Range("F40").Value = "Attracted"
$ws1.Range("F41").Value = "Impassioned"
$ws1.Range("F42").Value = "Impassioned"
$ws1.Range("F43").Value = "Attracted"
$ws1.Range("F44").Value = "Impassioned"
$ws1.Range("F46").Value = "Attracted"
$ws1.Range("F47").Value = "Impassioned"
$ws1.Range("F48").Value = "Impassioned"
$ws1.Range("F49").Value = "Impassioned"
$ws1.Range("F50").Value = "Impassioned"
$ws1.Range("F51").Value = "Impassioned"
$ws1.Range("F52").Value = "Impassioned"
$ws1.Range("F53").Value = "Impassioned"
$ws1.Range("F54").Value = "Attracted"
$ws1.Range("F55").Value = "Attracted"
$ws1.Range("F56").Value = "Impassioned"
$ws1.Range("F57").Value = "Impassioned"
$ws1.Range("F58").Value = "Impassioned"
$ws1.Range("F59").Value = "Apathetic"
$ws1.Range("F60").Value = "Attracted"
$ws1.Range("F61").Value = "Impassioned"
$ws1.Range("F62").Value = "Attracted"
$ws1.Range("F63").Value = "Impassioned"
$ws1.Range("F64").Value = "Impassioned"
$ws1.Range("F65").Value = "Impassioned"
$ws1.Range("F66").Value = "Impassioned"
$ws1.Range("F67").Value = "Attracted"
$ws1.Range("F68").Value = "Impassioned"
$ws1.Range("F69").Value = "Impassioned"
$ws1.Range("F70").Value = "Attracted"
$ws1.Range("F72").Value = "Impassioned"
$ws1.Range("F73").Value = "Impassioned"
$ws1.Range("F74").Value = "Attracted"
$ws1.Range("F75").Value = "Attracted"
$ws1.Range("F76").Value = "Attracted"
$ws1.Range("F79").Value = "Attracted"
$ws1.Range("F80").Value = "Impassioned"
$ws1.Range("F81").Value = "Attracted"
$ws1.Range("F82").Value = "Impassioned"
$ws1.Range("F83").Value = "Impassioned"
$ws1.Range("F84").Value = "Impassioned"
$ws1.Range("F85").Value = "Impassioned"
$ws1.Range("F86").Value = "Attracted"
$ws1.Range("F87").Value = "Apathetic"
$ws1.Range("F88").Value = "Apathetic"
$ws1.Range("F89").Value = "Impassioned"
$ws1.Range("F90").Value = "Attracted"
$ws1.Range("F91").Value = "Apathetic"
$ws1.Range("F92").Value = "Attracted"
$ws1.Range("F93").Value = "Impassioned"
$ws1.Range("F94").Value = "Attracted"
$ws1.Range("F95").Value = "Attracted"
$ws1.Range("F96").Value = "Impassioned"
$ws1.Range("F97").Value = "Impassioned"
$ws1.Range("F98").Value = "Impassioned"
$ws1.Range("F99").Value = "Impassioned"
$ws1.Range("F100").Value = "Impassioned"
$ws1.Range("F101").Value = "Impassioned"
$ws1.Range("F102").Value = "Impassioned"
$ws1.Range("F103").Value = "Impassioned"
$ws1.Range("F104").Value = "Impassioned"
$ws1.Range("F105").Value = "Impassioned"
$ws1.Range("F106").Value = "Impassioned"
$ws1.Range("F108").Value = "Impassioned"
$ws1.Range("F109").Value = "Impassioned"
$ws1.Range("F110").Value = "Apathetic"
$ws1.Range("F111").Value = "Apathetic"
$ws1.Range("F112").Value = "Apathetic"
$ws1.Range("F115").Value = "Impassioned"
$ws1.Range("F116").Value = "Impassioned"
$ws1.Range("F117").Value = "Impassioned"
$ws1.Range("F118").Value = "Apathetic"
$ws1.Range("F119").Value = "Impassioned"
$ws1.Range("F120").Value = "Impassioned"
$ws1.Range("F121").Value = "Impassioned"
$ws1.Range("F122").Value = "Impassioned"
$ws1.Range("F123").Value = "Impassioned"
$ws1.Range("F124").Value = "Impassioned"
$ws1.Range("F125").Value = "Impassioned"
$ws1.Range("F126").Value = "Impassioned"
$ws1.Range("F127").Value = "Impassioned"
$ws1.Range("F128").Value = "Impassioned"
$ws1.Range("F129").Value = "Impassioned"
$ws1.Range("F130").Value = "Impassioned"
$ws1.Range("F131").Value = "Impassioned"
$ws1.Range("F132").Value = "Attracted"
$ws1.Range("F133").Value = "Apathetic"
$ws1.Range("F134").Value = "Apathetic"
$ws1.Range("F135").Value = "Apathetic"
$ws1.Range("F136").Value = "Attracted"
$ws1.Range("F137").Value = "Impassioned"
$ws1.Range("F138").Value = "Attracted"
$ws1.Range("F139").Value = "Attracted"
$ws1.Range("F140").Value = "Apathetic"
$ws1.Range("F141").Value = "Impassioned"
$ws1.Range("F142").Value = "Impassioned"
$ws1.Range("F143").Value = "Impassioned"
$ws1.Range("F144").Value = "Impassioned"
$ws1.Range("F145").Value = "Attracted"
$ws1.Range("F146").Value = "Impassioned"
$ws1.Range("F147").Value = "Impassioned"
$ws1.Range("F148").Value = "Impassioned"
$ws1.Range("F149").Value = "Impassioned"
$ws1.Range("F150").Value = "Attracted"
$ws1.Range("F151").Value = "Impassioned"
$ws1.Range("F152").Value = "Impassioned"
$ws1.Range("F153").Value = "Attracted"
$ws1.Range("F154").Value = "Attracted"
$ws1.Range("F155").Value = "Impassioned"
$ws1.Range("F156").Value = "Impassioned"
$ws1.Range("F157").Value = "Impassioned"
$ws1.Range("F158").Value = "Impassioned"
$ws1.Range("F159").Value = "Apathetic"
$ws1.Range("F160").Value = "Apathetic"
$ws1.Range("F161").Value = "Impassioned"
$ws1.Range("F162").Value = "Impassioned"
$ws1.Range("F163").Value = "Apathetic"
$ws1.Range("F164").Value = "Impassioned"
$ws1.Range("F165").Value = "Impassioned"
$ws1.Range("F166").Value = "Impassioned"
$ws1.Range("F167").Value = "Impassioned"
$ws1.Range("F168").Value = "Impassioned"
$ws1.Range("F169").Value = "Impassioned"
$ws1.Range("F170").Value = "Apathetic"
$ws1.Range("F171").Value = "Impassioned"
$ws1.Range("F172").Value = "Attracted"
$ws1.Range("F173").Value = "Impassioned"
$ws1.Range("F174").Value = "Apathetic"
$ws1.Range("F175").Value = "Attracted"
$ws1.Range("F176").Value = "Impassioned"
$ws1.Range("F177").Value = "Attracted"
$ws1.Range("F178").Value = "Impassioned"
$ws1.Range("F179").Value = "Attracted"
$ws1.Range("F180").Value = "Attracted"
$ws1.Range("F181").Value = "Impassioned"
$ws1.Range("F182").Value = "Attracted"
$ws1.Range("F183").Value = "Impassioned"

# --- Sheet2 column F (Mindset/Test) data updates ---
$ws2.Range("F4").Value = "Attracted"
$ws2.Range("F5").Value = "Attracted"
$ws2.Range("F6").Value = "Impassioned"
$ws2.Range("F7").Value = "Attracted"
$ws2.Range("F13").Value = "Impassioned"
$ws2.Range("F15").Value = "Impassioned"
$ws2.Range("F17").Value = "Impassioned"
$ws2.Range("F19").Value = "Impassioned"
$ws2.Range("F20").Value = "Impassioned"
$ws2.Range("F26").Value = "Apathetic"
$ws2.Range("F27").Value = "Impassioned"
$ws2.Range("F28").Value = "Attracted"
$ws2.Range("F29").Value = "Attracted"
$ws2.Range("F30").Value = "Impassioned"
$ws2.Range("F31").Value = "Impassioned"
$ws2.Range("F32").Value = "Attracted"
$ws2.Range("F33").Value = "Impassioned"
$ws2.Range("F36").Value = "Attracted"
$ws2.Range("F37").Value = "Impassioned"
$ws2.Range("F38").Value = "Attracted"
$ws2.Range("F39").Value = "Impassioned"
$ws2.Range("F40").Value = "Attracted"
$ws2.Range("F41").Value = "Attracted"
$ws2.Range("F42").Value = "Apathetic"
$ws2.Range("F43").Value = "Attracted"
$ws2.Range("F44").Value = "Impassioned"
$ws2.Range("F45").Value = "Attracted"
$ws2.Range("F46").Value = "Impassioned"
$ws2.Range("F47").Value = "Impassioned"
$ws2.Range("F48").Value = "Attracted"
$ws2.Range("F49").Value = "Attracted"
$ws2.Range("F50").Value = "Apathetic"
$ws2.Range("F51").Value = "Apathetic"
$ws2.Range("F52").Value = "Apathetic"
$ws2.Range("F53").Value = "Impassioned"
$ws2.Range("F54").Value = "Impassioned"
$ws2.Range("F55").Value = "Impassioned"
$ws2.Range("F56").Value = "Impassioned"
$ws2.Range("F57").Value = "Attracted"
$ws2.Range("F58").Value = "Attracted"
$ws2.Range("F59").Value = "Attracted"
$ws2.Range("F60").Value = "Attracted"
$ws2.Range("F61").Value = "Attracted"
$ws2.Range("F62").Value = "Attracted"
$ws2.Range("F63").Value = "Impassioned"
$ws2.Range("F64").Value = "Impassioned"
$ws2.Range("F65").Value = "Impassioned"
$ws2.Range("F66").Value = "Attracted"
$ws2.Range("F67").Value = "Attracted"
$ws2.Range("F68").Value = "Impassioned"
$ws2.Range("F69").Value = "Impassioned"
$ws2.Range("F70").Value = "Impassioned"
$ws2.Range("F71").Value = "Impassioned"
$ws2.Range("F72").Value = "Impassioned"
$ws2.Range("F73").Value = "Impassioned"
$ws2.Range("F74").Value = "Impassioned"
$ws2.Range("F75").Value = "Impassioned"
$ws2.Range("F76").Value = "Attracted"
$ws2.Range("F77").Value = "Impassioned"
$ws2.Range("F78").Value = "Impassioned"
$ws2.Range("F79").Value = "Impassioned"
$ws2.Range("F80").Value = "Attracted"
$ws2.Range("F81").Value = "Impassioned"
$ws2.Range("F82").Value = "Impassioned"
$ws2.Range("F83").Value = "Impassioned"
$ws2.Range("F84").Value = "Impassioned"
$ws2.Range("F85").Value = "Impassioned"
$ws2.Range("F86").Value = "Impassioned"
$ws2.Range("F87").Value = "Impassioned"
$ws2.Range("F89").Value = "Impassioned"
$ws2.Range("F90").Value = "Impassioned"
$ws2.Range("F91").Value = "Impassioned"
$ws2.Range("F92").Value = "Attracted"
$ws2.Range("F93").Value = "Attracted"
$ws2.Range("F94").Value = "Attracted"
$ws2.Range("F97").Value = "Attracted"
$ws2.Range("F98").Value = "Impassioned"
$ws2.Range("F99").Value = "Attracted"
$ws2.Range("F101").Value = "Impassioned"
$ws2.Range("F102").Value = "Impassioned"
$ws2.Range("F103").Value = "Impassioned"
$ws2.Range("F104").Value = "Attracted"
$ws2.Range("F105").Value = "Impassioned"
$ws2.Range("F106").Value = "Impassioned"
$ws2.Range("F107").Value = "Impassioned"
$ws2.Range("F108").Value = "Apathetic"
$ws2.Range("F110").Value = "Apathetic"
$ws2.Range("F111").Value = "Impassioned"
$ws2.Range("F112").Value = "Attracted"
$ws2.Range("F113").Value = "Apathetic"
$ws2.Range("F115").Value = "Attracted"
$ws2.Range("F116").Value = "Apathetic"
$ws2.Range("F117").Value = "Impassioned"
$ws2.Range("F118").Value = "Impassioned"
$ws2.Range("F119").Value = "Impassioned"
$ws2.Range("F120").Value = "Apathetic"
$ws2.Range("F121").Value = "Apathetic"
$ws2.Range("F122").Value = "Attracted"
$ws2.Range("F123").Value = "Attracted"
$ws2.Range("F124").Value = "Impassioned"
$ws2.Range("F125").Value = "Impassioned"
$ws2.Range("F126").Value = "Impassioned"
$ws2.Range("F127").Value = "Impassioned"
$ws2.Range("F128").Value = "Apathetic"
$ws2.Range("F129").Value = "Apathetic"
$ws2.Range("F130").Value = "Impassioned"
$ws2.Range("F131").Value = "Impassioned"
$ws2.Range("F132").Value = "Impassioned"
$ws2.Range("F133").Value = "Impassioned"
$ws2.Range("F134").Value = "Impassioned"
$ws2.Range("F135").Value = "Impassioned"
$ws2.Range("F136").Value = "Impassioned"
$ws2.Range("F137").Value = "Impassioned"
$ws2.Range("F138").Value = "Impassioned"
$ws2.Range("F141").Value = "Apathetic"
$ws2.Range("F143").Value = "Apathetic"
$ws2.Range("F144").Value = "Apathetic"
$ws2.Range("F145").Value = "Apathetic"
$ws2.Range("F146").Value = "Impassioned"
$ws2.Range("F147").Value = "Impassioned"
$ws2.Range("F148").Value = "Impassioned"
$ws2.Range("F149").Value = "Apathetic"
$ws2.Range("F150").Value = "Apathetic"
$ws2.Range("F151").Value = "Impassioned"
$ws2.Range("F152").Value = "Impassioned"
$ws2.Range("F153").Value = "Impassioned"
$ws2.Range("F154").Value = "Impassioned"
$ws2.Range("F155").Value = "Impassioned"
$ws2.Range("F156").Value = "Attracted"
$ws2.Range("F157").Value = "Attracted"
$ws2.Range("F158").Value = "Impassioned"
$ws2.Range("F159").Value = "Impassioned"
$ws2.Range("F160").Value = "Impassioned"
$ws2.Range("F161").Value = "Impassioned"
$ws2.Range("F162").Value = "Apathetic"
$ws2.Range("F163").Value = "Attracted"
$ws2.Range("F164").Value = "Impassioned"
$ws2.Range("F165").Value = "Impassioned"
$ws2.Range("F166").Value = "Apathetic"
$ws2.Range("F167").Value = "Impassioned"
$ws2.Range("F168").Value = "Impassioned"
$ws2.Range("F169").Value = "Impassioned"
$ws2.Range("F170").Value = "Impassioned"
$ws2.Range("F171").Value = "Impassioned"
$ws2.Range("F172").Value = "Impassioned"
$ws2.Range("F173").Value = "Impassioned"
$ws2.Range("F174").Value = "Attracted"
$ws2.Range("F175").Value = "Impassioned"
$ws2.Range("F176").Value = "Impassioned"
$ws2.Range("F177").Value = "Impassioned"
$ws2.Range("F178").Value = "Apathetic"
$ws2.Range("F179").Value = "Impassioned"
$ws2.Range("F180").Value = "Apathetic"
$ws2.Range("F181").Value = "Apathetic"
$ws2.Range("F182").Value = "Attracted"
$ws2.Range("F183").Value = "Impassioned"
